# RPA datasets push 2024-07-05
# The "02_38커뮤니케이션(최근일자기준)" sheet lists IPO book-building
# entries ordered by 수요예측일 (book-building date). The 뱅크웨어글로벌
# (BankWare Global) entry had its book-building date corrected from
# "2024.07.08~07.12" to "2024.07.23~07.29", which moves the row up in
# the date-ordered list: from old row 13 to new row 5 (just above
# 티디에스팜), shifting the rows in between down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a fresh blank row at row 5, pushing the existing rows 5..21 down to 6..22.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the (corrected) BankWare Global entry.
$ws.Cells.Item(5, 1).Value = "뱅크웨어글로벌"
$ws.Cells.Item(5, 2).Value = "2024.07.23~07.29"
$ws.Cells.Item(5, 3).Value = "16,000~19,000"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = 22400
$ws.Cells.Item(5, 6).Value = "미래에셋증권"

# The original BankWare Global row (previously row 13) is now row 14 after
# the insert above; remove it so the entry only appears once, in its new
# position, restoring the sheet back to 21 data rows.
$ws.Rows.Item(14).Delete()
